$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.0000000000000000002330246048421937
$ws.Range("C3").Value = 0.4458327234947346
$ws.Range("C4").Value = 0.06357623681831188
$ws.Range("C5").Value = 0.000000000000000008807695185509702
$ws.Range("C6").Value = 0.3486403815655779
$ws.Range("C7").Value = 0.0000000000000000003661815218948759
$ws.Range("C8").Value = 0.0000000000000000009653876486319454
$ws.Range("C9").Value = 0.05261526521361033
$ws.Range("C10").Value = 0.02203726983585222
$ws.Range("C11").Value = 0.05890727244953593
$ws.Range("C13").Value = 0.008390850622377223
